# CSAT Performance Reports - add new day/agent data
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: Daywise_Report  (A1:F8 -> A1:F10)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Daywise_Report")

# Update the MTD summary row (row 2)
$ws1.Range("B2").Value = 8
$ws1.Range("D2").Value = 9
$ws1.Range("E2").Value = 19
$ws1.Range("F2").Value = 2.05

# New row 9 (2025-10-08)
$ws1.Range("A3").Copy()
$ws1.Range("A9").PasteSpecial(-4122)
$ws1.Range("A9").Value = 45938

$ws1.Range("D3").Copy()
$ws1.Range("B9").PasteSpecial(-4122)
$ws1.Range("B9").Value = 1

$ws1.Range("C3").Copy()
$ws1.Range("C9").PasteSpecial(-4122)
$ws1.Range("C9").Value = 0

$ws1.Range("C3").Copy()
$ws1.Range("D9").PasteSpecial(-4122)
$ws1.Range("D9").Value = 0

$ws1.Range("C3").Copy()
$ws1.Range("E9").PasteSpecial(-4122)
$ws1.Range("E9").Value = 1

$ws1.Range("D3").Copy()
$ws1.Range("F9").PasteSpecial(-4122)
$ws1.Range("F9").Value = 1

# New row 10 (2025-10-09)
$ws1.Range("A3").Copy()
$ws1.Range("A10").PasteSpecial(-4122)
$ws1.Range("A10").Value = 45939

$ws1.Range("C3").Copy()
$ws1.Range("B10").PasteSpecial(-4122)
$ws1.Range("B10").Value = 0

$ws1.Range("C3").Copy()
$ws1.Range("C10").PasteSpecial(-4122)
$ws1.Range("C10").Value = 0

$ws1.Range("C2").Copy()
$ws1.Range("D10").PasteSpecial(-4122)
$ws1.Range("D10").Value = 2

$ws1.Range("C3").Copy()
$ws1.Range("E10").PasteSpecial(-4122)
$ws1.Range("E10").Value = 2

$ws1.Range("F2").Copy()
$ws1.Range("F10").PasteSpecial(-4122)
$ws1.Range("F10").Value = 3

$ws1.ListObjects.Item(1).Resize($ws1.Range("A1:F10"))

Write-Host "Sheet1 done"

# ---------------------------------------------------------------------------
# Sheet 2: Agentwise_Report  (A1:F16 -> A1:F18)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Agentwise_Report")
$lo2 = $ws2.ListObjects.Item(1)

# Update row 3 (SBM202): B3 0->1 (gains red fill), E3 1->2, F3 3->2 (stays green)
$ws2.Range("D3").Copy()
$ws2.Range("B3").PasteSpecial(-4122)
$ws2.Range("B3").Value = 1
$ws2.Range("E3").Value = 2
$ws2.Range("F3").Value = 2

# Insert new row for agent VPS214 before current row 12 (VPS228), shifting rows 12-16 down to 13-17
$ws2.Rows("12:12").Insert()
$ws2.Range("A12").Value = "VPS214"
$ws2.Range("B12").Value = 0
$ws2.Range("C12").Value = 0
$ws2.Range("D12").Value = 1
$ws2.Range("E12").Value = 1
$ws2.Range("F12").Value = 3

# Append new row for agent VPS264 at the end (row 18)
# copy formatting for row18 from row12 (same B/C/D/E/F fill pattern), then set values
$ws2.Range("A12:F12").Copy()
$ws2.Range("A18:F18").PasteSpecial(-4122)
$ws2.Range("A18").Value = "VPS264"
$ws2.Range("B18").Value = 0
$ws2.Range("C18").Value = 0
$ws2.Range("D18").Value = 1
$ws2.Range("E18").Value = 1
$ws2.Range("F18").Value = 3

$lo2.Resize($ws2.Range("A1:F18"))

Write-Host "Sheet2 done"

# ---------------------------------------------------------------------------
# Sheet 3: Negative_Responses  (A1:R8 -> A1:R9)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Negative_Responses")
$lo3 = $ws3.ListObjects.Item(1)

$ws3.Range("A9").Value = "J082665608785"
$ws3.Range("B9").Value = 9445345417
$ws3.Range("C9").Value = "08-10-2025 07:14 PM"
$ws3.Range("D9").Value = "CSAT 1"
$ws3.Range("E9").Value = "08-10-2025 07:19 PM"
$ws3.Range("H9").Value = "QUERY"
$ws3.Range("I9").Value = "PREPAID"
$ws3.Range("J9").Value = "SERVICES"
$ws3.Range("K9").Value = "PLANS AND VALIDITY DETAILS"
$ws3.Range("L9").Value = "ENQUIRY ON BALANCE AND VALIDITY"
$ws3.Range("M9").Value = "BSS_WEST OA_1"
$ws3.Range("N9").Value = "CHENNAI"
$ws3.Range("O9").Value = "TAMILNADU"
$ws3.Range("P9").Value = "SBM202"
$ws3.Range("Q9").Value = 1

$ws3.Range("R8").Copy()
$ws3.Range("R9").PasteSpecial(-4122)
$ws3.Range("R9").Value = 45938

$lo3.Resize($ws3.Range("A1:R9"))

Write-Host "Sheet3 done"

# ---------------------------------------------------------------------------
# Sheet 4: Daywise_Agent_Performance  (A1:G17 -> A1:G20)
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Daywise_Agent_Performance")
$lo4 = $ws4.ListObjects.Item(1)

# Row 18 (2025-10-08, SBM202) - same style pattern as row 5 (A=date,C=red,G=red)
$ws4.Range("A5:G5").Copy()
$ws4.Range("A18:G18").PasteSpecial(-4122)
$ws4.Range("A18").Value = 45938
$ws4.Range("B18").Value = "SBM202"
$ws4.Range("C18").Value = 1
$ws4.Range("D18").Value = 0
$ws4.Range("E18").Value = 0
$ws4.Range("F18").Value = 1
$ws4.Range("G18").Value = 1

# Row 19 (2025-10-09, VPS214) - same style pattern as row 2 (A=date,G=green)
$ws4.Range("A2:G2").Copy()
$ws4.Range("A19:G19").PasteSpecial(-4122)
$ws4.Range("A19").Value = 45939
$ws4.Range("B19").Value = "VPS214"
$ws4.Range("C19").Value = 0
$ws4.Range("D19").Value = 0
$ws4.Range("E19").Value = 1
$ws4.Range("F19").Value = 1
$ws4.Range("G19").Value = 3

# Row 20 (2025-10-09, VPS264) - same style pattern as row 2 (A=date,G=green)
$ws4.Range("A2:G2").Copy()
$ws4.Range("A20:G20").PasteSpecial(-4122)
$ws4.Range("A20").Value = 45939
$ws4.Range("B20").Value = "VPS264"
$ws4.Range("C20").Value = 0
$ws4.Range("D20").Value = 0
$ws4.Range("E20").Value = 1
$ws4.Range("F20").Value = 1
$ws4.Range("G20").Value = 3

$lo4.Resize($ws4.Range("A1:G20"))

Write-Host "Sheet4 done"
